# Add a new "credits" textbox (TextBox 13) to slide 1, matching the
# shape that was uploaded in the target commit.
#
# NOTE on units: this COM-interop runtime's Shapes.AddTextbox / shape
# Left-Top-Width-Height style setters take values in *points* and store
# them internally as EMU (1 pt = 12700 EMU). The target OOXML gives EMU
# offsets directly, so we convert EMU -> points (divide by 12700) before
# calling into the COM layer.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700.0

$offX = 9540512 / $emuPerPt
$offY = 453207 / $emuPerPt
$extCX = 2533271 / $emuPerPt
$extCY = 577081 / $emuPerPt

# The deck already has shapes with ids 1,2,3,5,7,8,9 (ids 4 and 6 are
# "holes" left by shapes deleted earlier in the deck's history). The
# target shape must land on id=10, so first consume the two holes with
# throwaway textboxes, delete them, and only then add the real shape -
# the id allocator fills holes before minting a fresh id.
$filler1 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$filler2 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$filler1.Delete()
$filler2.Delete()

$tb = $s.Shapes.AddTextbox(1, $offX, $offY, $extCX, $extCY)
$tb.Name = "TextBox 13"

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

# noFill
$tb.Fill.Visible = $false

# Setting multi-paragraph text in one shot (text containing `r) only
# stamps the default lang="en-US" onto the first paragraph's run, so
# build it up paragraph-by-paragraph and restamp Size/LanguageID on
# each newly-inserted span to keep every run's formatting consistent.
$tr = $tf.TextRange
$tr.Text = "Ansan Technical High School"
$tr.Font.Size = 10.5
$tr.LanguageID = "en-US"

$r2 = $tf.TextRange.InsertAfter("`rDept. Computer")
$r2.Font.Size = 10.5
$r2.LanguageID = "en-US"

$r3 = $tf.TextRange.InsertAfter("`rMade by kig2929kig@gmail.com")
$r3.Font.Size = 10.5
$r3.LanguageID = "en-US"

Write-Host "Shape count:" $s.Shapes.Count
Write-Host "New shape id:" $tb.Id
Write-Host "New shape name:" $tb.Name
